$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the missing time-log entry for row 24 (date 42971 -> 4 hours worked).
# F3's SUM(D3:D33) formula picks this up automatically on recalculation
# (94.5 -> 98.5).
$ws.Range("D24").Value = 4

# Scroll the view down so row 10 is the top-left visible row, and move the
# active selection/cursor to H22.
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H22").Select()
